$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "short-url" column (B) shared-string content change: ve0E4M -> P1VV4y ---
# Every data row (2-368) shares this same text, so rewrite the whole column at
# once so it collapses back onto a single shared string (matches the source
# edit, which just edited the shared string's text in place).
$ws.Range("B2:B368").Value = "P1VV4y"

# --- Numeric-looking text cells that must stay text (shared-string) cells ---
# Plain `.Value = "32"` gets auto-typed as a number by Excel, so force the
# cell to Text format first, then assign the digit string.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "N361" "32"
Set-TextValue "O363" "9"
Set-TextValue "N366" "281"
Set-TextValue "O366" "3576"
Set-TextValue "N368" "24"
Set-TextValue "O368" "245"
